$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update simulation result values for the 380 kV case (Case_0_186, pl_mw.xlsx)
$data = @{
    "C2"=0.04851133617081871; "D2"=0.1721232221067623; "E2"=0.2585910595507244; "F2"=1.566472636612914; "G2"=1.283304412812583; "H2"=1.165733780664482; "J2"=0.4030933344628522; "M2"=11.03267308633696;
    "C3"=0.04303509047016973; "D3"=0.1559118863778508; "E3"=0.2326807926679066; "F3"=1.629441323364802; "G3"=1.306603723938878; "H3"=1.196388683199899; "J3"=0.3609169052584775; "M3"=9.728459456845826;
    "C4"=0.0397026064644308; "D4"=0.1459268339486357; "E4"=0.2168485604219512; "F4"=1.672153050430197; "G4"=1.324581164255704; "H4"=1.217492531298589; "J4"=0.3352721812424022; "M4"=8.925665515002549;
    "C5"=0.03835168176055959; "D5"=0.1418497945449531; "E5"=0.2104146866927152; "F5"=1.690554156265343; "G5"=1.332806048665617; "H5"=1.226655141088372; "J5"=0.3248808200953022; "M5"=8.597962056204153;
    "C6"=0.03812777728740002; "D6"=0.1411723140413557; "E6"=0.2093473886255524; "F6"=1.693669140916548; "G6"=1.334225376550705; "H6"=1.228210245793946; "J6"=0.3231587924318831; "M6"=8.543511915694978;
    "C7"=0.03968435929022007; "D7"=0.1458718823950989; "E7"=0.216761720287991; "F7"=1.672397211671139; "G7"=1.324688479918507; "H7"=1.21761383793104; "J7"=0.335131805904183; "M7"=8.921248321412122;
    "C8"=0.04661663758474788; "D8"=0.166539899943416; "E8"=0.249640360704646; "F8"=1.587330345972347; "G8"=1.290561187468825; "H8"=1.175823287887681; "J8"=0.3884959844843081; "M8"=10.58336038017137;
    "C9"=0.06046976379906255; "D9"=0.2068347990554287; "E9"=0.3147988307733982; "F9"=1.453664025653268; "G9"=1.253889171928535; "H9"=1.112485703623321; "J9"=0.4953531201156807; "M9"=13.82979841242292;
    "C10"=0.07083702209192211; "D10"=0.2363194942805933; "E10"=0.3632062189319925; "F10"=1.37711790708299; "G10"=1.246977042139093; "H10"=1.078020812397256; "J10"=0.5755349561040646; "M10"=16.21175574969953;
    "C11"=0.07560173214649524; "D11"=0.2497128687702457; "E11"=0.3853737728053943; "F11"=1.347325671664478; "G11"=1.248541243293261; "H11"=1.065125630745854; "J11"=0.6124544521571806; "M11"=17.29584229795364;
    "C12"=0.07741359766639277; "D12"=0.2547823039495256; "E12"=0.3937917456983229; "F12"=1.336795643063525; "G12"=1.249841577108924; "H12"=1.060656844011703; "J12"=0.6265057809404482; "M12"=17.70653411182411;
    "C13"=0.07702303357626761; "D13"=0.2536906071664475; "E13"=0.3919776929764112; "F13"=1.339029578257112; "G13"=1.249529533231652; "H13"=1.061600618520743; "J13"=0.6234763168105815; "M13"=17.61807501515085;
    "C14"=0.0757506401824628; "D14"=0.2501299785605227; "E14"=0.3860658364558844; "F14"=1.346444128793522; "G14"=1.248633851767181; "H14"=1.064749594154165; "J14"=0.6136090053196597; "M14"=17.32962613495027;
    "C15"=0.07497226658675515; "D15"=0.2479486991958311; "E15"=0.3824478048408508; "F15"=1.351084505320728; "G15"=1.248178361914029; "H15"=1.066732820430929; "J15"=0.6075744124210587; "M15"=17.15296824538353;
    "C16"=0.07052666435421884; "D16"=0.2354438525911178; "E16"=0.3617606792694517; "F16"=1.379168505868222; "G16"=1.246972182696368; "H16"=1.078920741722357; "J16"=0.5731316935072357; "M16"=16.14092673091295;
    "C17"=0.06781228552564755; "D17"=0.2277678725134251; "E17"=0.3491091020129033; "F17"=1.397705100927595; "G17"=1.247459840801469; "H17"=1.087120340669173; "J17"=0.5521207097835372; "M17"=15.5202720179671;
    "C18"=0.06625558374896912; "D18"=0.2233510001934462; "E18"=0.3418459041340185; "F18"=1.408838444556579; "G18"=1.248181603445772; "H18"=1.09209761243639; "J18"=0.5400771368043991; "M18"=15.16333071957922;
    "C19"=0.06572927311189858; "D19"=0.2218551925134875; "E19"=0.3393889828034276; "F19"=1.412688148094205; "G19"=1.248500846034034; "H19"=1.093827249583171; "J19"=0.5360062989711309; "M19"=15.04248107993811;
    "C20"=0.06810076204962456; "D20"=0.228585182303334; "E20"=0.3504544534059875; "F20"=1.395682838666929; "G20"=1.24736203827905; "H20"=1.086220357933627; "J20"=0.5543530334618367; "M20"=15.58633665968165;
    "C21"=0.07612416263155808; "D21"=0.2511758813990639; "E21"=0.3878016305683616; "F21"=1.34424566367467; "G21"=1.248877467763094; "H21"=1.063813303293131; "J21"=0.6165052978428207; "M21"=17.41434509058314;
    "C22"=0.0814123523916237; "D22"=0.2659267790059232; "E22"=0.4123490465951392; "F22"=1.315023674077338; "G22"=1.254009895899969; "H22"=1.051591033250702; "J22"=0.6575411144887937; "M22"=18.61010622952659;
    "C23"=0.07858567979705811; "D23"=0.2580550184335948; "E23"=0.3992340420890912; "F23"=1.330208140040057; "G23"=1.250881099628941; "H23"=1.057887859734791; "J23"=0.6355991044486018; "M23"=17.97177566256789;
    "C24"=0.06797032992503205; "D24"=0.2282156885683548; "E24"=0.3498461879505186; "F24"=1.396595621672191; "G24"=1.247404882257854; "H24"=1.086626421341606; "J24"=0.5533436885490062; "M24"=15.55646921189793;
    "C25"=0.05669122655658043; "D25"=0.1959568747581102; "E25"=0.2970869520061115; "F25"=1.486128561463389; "G25"=1.260410334037488; "H25"=1.127554328933343; "J25"=0.4661766211201552; "M25"=12.95247859135213;
}

foreach ($key in $data.Keys) {
    $ws.Range($key).Value = $data[$key]
}

Write-Host "Updated $($data.Count) cells"